# GROUP9_BOM.xlsx update
# The "Voltage booster" line item (row 19) count drops from 2 to 1.
# This is the single substantive edit in the commit; Excel will
# auto-recalculate the dependent formulas (F19 line total, F21 the
# POWER section subtotal, and F40 the grand total).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Voltage booster count: 2 -> 1
$ws.Range("D19").Value = 1

# Reflect the cursor/selection position left behind by the author when
# the workbook was last saved.
$ws.Range("H20").Select()
